$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# EMU -> point conversion factor used throughout the PowerPoint object model
$emuPerPt = 12700

# 1) Move the results picture up a bit (x/size unchanged)
$pic = $s.Shapes.Item(2)
$pic.Top = 2785685 / $emuPerPt

# 2) Add the first new textbox: euclidean-distance note (centered, word-wrapped)
$tb1 = $s.Shapes.AddTextbox(1, 1758950 / $emuPerPt, 5207000 / $emuPerPt, 8674100 / $emuPerPt, 369332 / $emuPerPt)
$tb1.Fill.Visible = 0
$tb1.TextFrame.WordWrap = -1
$tb1.TextFrame.AutoSize = 1

$tb1r = $tb1.TextFrame.TextRange
$tb1r.Text = "With"
$tb1r.LanguageID = "fr-FR"
$tb1r.ParagraphFormat.Alignment = 2

$tb1r2 = $tb1.TextFrame.TextRange.InsertAfter(" a simple euclidien distance : (0.45614035087719296, 0.38461538461538464)")
$tb1r2.LanguageID = "fr-FR"

# 3) Add the second new textbox: best-K note (no wrap)
$tb2 = $s.Shapes.AddTextbox(1, 4047843 / $emuPerPt, 5576332 / $emuPerPt, 4096314 / $emuPerPt, 369332 / $emuPerPt)
$tb2.Fill.Visible = 0
$tb2.TextFrame.WordWrap = 0
$tb2.TextFrame.AutoSize = 1

$tb2r = $tb2.TextFrame.TextRange
$tb2r.Text = "Best K "
$tb2r.LanguageID = "fr-FR"

$tb2r2 = $tb2.TextFrame.TextRange.InsertAfter("found")
$tb2r2.LanguageID = "fr-FR"

$tb2r3 = $tb2.TextFrame.TextRange.InsertAfter(" ")
$tb2r3.LanguageID = "fr-FR"

$tb2r4 = $tb2.TextFrame.TextRange.InsertAfter("with")
$tb2r4.LanguageID = "fr-FR"

$tb2r5 = $tb2.TextFrame.TextRange.InsertAfter(" a 1000 values set : 100")
$tb2r5.LanguageID = "fr-FR"
